# Atualizacao de bases das ligas, do dia: 07-03-2024 as 23:43
# Adds results/odds for finished matches on rows 123-125 (home/away goals,
# result letter, and updated closing odds) and appends 5 new upcoming
# fixtures as rows 126-130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue {
    param($Sheet, [int]$Row, [int]$Col, $Value)
    if ($Value -ne $null) {
        $Sheet.Cells.Item($Row, $Col).Value = $Value
    }
}

function Write-DataRow {
    param($Sheet, [int]$Row, $Data)

    # Columns A-G are assumed to already exist (id/meta columns) for rows
    # that were already present before this edit; for brand-new rows the
    # caller is responsible for seeding A/E styles beforehand.
    Set-CellValue $Sheet $Row 8  $Data.H
    Set-CellValue $Sheet $Row 9  $Data.I
    Set-CellValue $Sheet $Row 10 $Data.J
    Set-CellValue $Sheet $Row 11 $Data.K
    Set-CellValue $Sheet $Row 12 $Data.L
    Set-CellValue $Sheet $Row 13 $Data.M
    Set-CellValue $Sheet $Row 14 $Data.N
    Set-CellValue $Sheet $Row 15 $Data.O
    Set-CellValue $Sheet $Row 16 $Data.P
    Set-CellValue $Sheet $Row 17 $Data.Q
    Set-CellValue $Sheet $Row 18 $Data.R
    Set-CellValue $Sheet $Row 19 $Data.S
    Set-CellValue $Sheet $Row 20 $Data.T
    Set-CellValue $Sheet $Row 21 $Data.U
    Set-CellValue $Sheet $Row 22 $Data.V
    Set-CellValue $Sheet $Row 23 $Data.W
    Set-CellValue $Sheet $Row 24 $Data.X
    Set-CellValue $Sheet $Row 25 $Data.Y
    Set-CellValue $Sheet $Row 26 $Data.Z
    Set-CellValue $Sheet $Row 27 $Data.AA
    Set-CellValue $Sheet $Row 28 $Data.AB
    Set-CellValue $Sheet $Row 29 $Data.AC
}

# ---------------------------------------------------------------------------
# 1) Rows that already existed (123-125): fill in the final score (H/I),
#    result letter (J) and refresh the odds columns (K-AC) now that the
#    matches have been played.
# ---------------------------------------------------------------------------

Write-DataRow $ws 123 @{
    H=2; I=0; J="H"
    K=1.25; L=4.75; M=11
    N=1.2; O=5.5; P=13
    Q=-1.75
    R=1.95; S=1.9
    T=2.75
    U=1.975; V=1.875
    W=0.2; X=-1; Y=-1; Z=0.475; AA=-0.5; AB=-1; AC=0.875
}

Write-DataRow $ws 124 @{
    H=1; I=1; J="D"
    K=3.6; L=3.3; M=1.909
    N=4.5; O=3.2; P=1.8
    Q=0.5
    R=1.975; S=1.875
    T=2.25
    U=1.975; V=1.875
    W=-1; X=2.2; Y=-1; Z=0.9750000000000001; AA=-1; AB=-0.5; AC=0.4375
}

Write-DataRow $ws 125 @{
    H=2; I=1; J="H"
    K=3.3; L=3.2; M=2.1
    N=4.333; O=3.2; P=1.85
    Q=0.5
    R=1.925; S=1.925
    T=2.25
    U=2; V=1.85
    W=3.333; X=-1; Y=-1; Z=0.925; AA=-1; AB=1; AC=-1
}

# ---------------------------------------------------------------------------
# 2) Brand-new upcoming fixtures, appended as rows 126-130. Column A needs
#    the bold/centered/bordered "id" style (same as used through row 125)
#    and column E needs the date-time number format style; copy both from
#    row 125 so the new cells pick up the exact same style indices instead
#    of minting new ones.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row=126; Id=124; MatchId=6788929; Date=45359.54166666666; Home="NK Varazdin";          Away="Istra 1961";
       K=1.909; L=3.4;  M=3.6;  N=2.2;  O=3.2;  P=3.1;  Q=-0.25; R=1.975; S=1.875; T=2.25; U=2;     V=1.85;  W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=127; Id=125; MatchId=6769304; Date=45360.45833333334; Home="NK Lokomotiva Zagreb"; Away="NK Rudes";
       K=1.333; L=4.75; M=7.5;  N=1.3;  O=4.5;  P=9.5;  Q=-1.5;  R=2;     S=1.85;  T=2.75; U=2.025; V=1.825; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=128; Id=126; MatchId=6788930; Date=45360.54861111111; Home="HNK Gorica";            Away="Hajduk Split";
       K=6;     L=3.75; M=1.533; N=7;    O=3.75; P=1.45; Q=1;     R=1.95;  S=1.9;   T=2.25; U=1.875; V=1.975; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=129; Id=127; MatchId=6788932; Date=45361.45833333334; Home="HNK Rijeka";            Away="NK Osijek";
       K=1.533; L=3.75; M=6;    N=1.55; O=3.8;  P=5.75; Q=-1;    R=2;     S=1.85;  T=2.5;  U=1.925; V=1.925; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ Row=130; Id=128; MatchId=6788931; Date=45361.5625;         Home="Dinamo Zagreb";         Away="Slaven Belupo";
       K=1.25;  L=6;    M=8;    N=1.222;O=6;    P=9.5;  Q=-1.75; R=1.875; S=1.975; T=3;    U=1.975; V=1.875; W=0; X=0; Y=0; Z=0; AA=0 }
)

foreach ($nr in $newRows) {
    $row = $nr.Row

    # Seed column A (id, bold/center/border style) and column E (date style)
    # by copying the formatting from row 125, which already carries the
    # right style indices; this avoids creating brand-new style entries.
    $ws.Cells.Item(125, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item(125, 5).Copy($ws.Cells.Item($row, 5))

    $ws.Cells.Item($row, 1).Value = $nr.Id
    $ws.Cells.Item($row, 2).Value = $nr.MatchId
    $ws.Cells.Item($row, 3).Value = "Croatia HNL"
    $ws.Cells.Item($row, 4).Value = "Croatia HNL"
    $ws.Cells.Item($row, 5).Value = $nr.Date
    $ws.Cells.Item($row, 6).Value = $nr.Home
    $ws.Cells.Item($row, 7).Value = $nr.Away

    Write-DataRow $ws $row $nr
}
